$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-04 18:56:38"

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
